$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.384.11'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '3.908.73'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '527.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +9.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.52'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.615'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.24%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('E10').Value = '  -4.85%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000338'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '42.07'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.53%  '
$ws.Range('D13').Value = '4.532.73'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.26'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').Value = '3.918.65'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('E16').Value = '  +9.58%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.01'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '19.76'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('D20').Value = '69.275.92'
$ws.Range('E20').Value = '  +1.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '427.35'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('E22').Value = '  -5.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '88.62'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').Value = '  -4.60%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.04'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +10.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.45'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -6.54%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.62'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.85%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.48'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.35%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '678.46'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.67%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '13.14'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.55%  '
$ws.Range('E31').Value = '  -3.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.83'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '68.84'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +11.40%  '
$ws.Range('D34').Value = '0.0₃0891'
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.434'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +9.63%  '
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.149'
$ws.Range('D38').Style = "Normal"
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.23'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.21%  '
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('E43').Value = '  +6.90%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.80'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -6.93%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.000292'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +20.94%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.36'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0352'
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.98'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +6.54%  '
$ws.Range('D50').Value = '2.745.85'
$ws.Range('E50').Value = '  +13.91%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '145.15'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.27%  '
